$wb = $excel.ActiveWorkbook

# "Sheet2" (the 4th sheet / xl/worksheets/sheet4.xml) becomes the active/selected
# sheet; this also clears tabSelected from whichever sheet was previously active
# ("Sheet1" / xl/worksheets/sheet3.xml) and updates workbook.xml's activeTab.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

# Append a second "table block" below the existing one - a st_level meta/header
# row followed by a header row and 5 data rows. The meta row (row 11) only gets
# a table name + table-name-id (no JSON meta info in column C), matching
# "support empty meta info".
$ws2.Range("A11").Value = "st_level"
$ws2.Range("B11").Value = "st_levelTable"

$ws2.Range("A12").Value = "id"
$ws2.Range("B12").Value = "nameenum:Kmyenum{A,B=21,C}"
$ws2.Range("C12").Value = "my:name"

$ws2.Range("A13").Value = 0
$ws2.Range("B13").Value = "A"
$ws2.Range("C13").Value = "D=4"

$ws2.Range("A14").Value = 1
$ws2.Range("B14").Value = "A"
$ws2.Range("C14").Value = "E"

$ws2.Range("A15").Value = 2
$ws2.Range("B15").Value = "A"
$ws2.Range("C15").Value = "F"

$ws2.Range("A16").Value = 3
$ws2.Range("B16").Value = "A"
$ws2.Range("C16").Value = "D"

$ws2.Range("A17").Value = 4
$ws2.Range("B17").Value = "A"
$ws2.Range("C17").Value = "E"

# Reflect the new selection/highlighted block on the sheet.
$ws2.Range("A11:C17").Select()
